{"js": "// Remove the unnecessary word \"the\" from \"far from the campus\" ->\n// \"far from campus\" in the Collaboration Challenges paragraph.\nconst body = context.document.body;\n\nconst results = body.search(\"from the campus\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target phrase \"from the campus\" not found.');\n}\n\nresults.items[0].insertText(\"from campus\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Removed unnecessary word from collab challenge:\n# \"far from the campus\" -> \"far from campus\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"from the campus\",  # FindText\n    $false,             # MatchCase\n    $false,             # MatchWholeWord\n    $false,             # MatchWildcards\n    $false,             # MatchSoundsLike\n    $false,             # MatchAllWordForms\n    $true,               # Forward\n    1,                   # Wrap (wdFindContinue)\n    $false,              # Format\n    \"from campus\",       # ReplaceWith\n    2                    # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    throw 'Target phrase \"from the campus\" not found.'\n}\n"}
